$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-12 22:48:45"
$ws.Range("E3").Value = "2026-02-12 22:48:47"
$ws.Range("E4").Value = "2026-02-12 22:48:50"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "37%"
$ws.Range("J4").Value = "999.9 hPa"
$ws.Range("E5").Value = "2026-02-12 22:48:53"
$ws.Range("E6").Value = "2026-02-12 22:48:55"
$ws.Range("J6").Value = "999.7 hPa"
$ws.Range("E7").Value = "2026-02-12 22:48:58"
$ws.Range("E8").Value = "2026-02-12 22:49:01"
$ws.Range("J8").Value = "1001.8 hPa"
$ws.Range("E9").Value = "2026-02-12 22:49:04"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "64%"
$ws.Range("O9").Value = "12.8 °C"
$ws.Range("E10").Value = "2026-02-12 22:49:06"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "48%"
$ws.Range("O10").Value = "14.6 °C"
$ws.Range("E11").Value = "2026-02-12 22:49:09"
$ws.Range("O11").Value = "9.1 °C"
$ws.Range("E12").Value = "2026-02-12 22:49:12"
$ws.Range("O12").Value = "12.4 °C"
$ws.Range("E13").Value = "2026-02-12 22:49:14"
$ws.Range("J13").Value = "1002.4 hPa"
$ws.Range("N13").Value = "4.0 °C 22:14 TU"
$ws.Range("E14").Value = "2026-02-12 22:49:16"
$ws.Range("E15").Value = "2026-02-12 22:49:19"
$ws.Range("E16").Value = "2026-02-12 22:49:22"
$ws.Range("O16").Value = "-4.5 °C"
$ws.Range("E17").Value = "2026-02-12 22:49:24"
$ws.Range("E18").Value = "2026-02-12 22:49:27"
$ws.Range("J18").Value = "1000.1 hPa"
$ws.Range("O18").Value = "16.4 °C"
$ws.Range("E19").Value = "2026-02-12 22:49:30"
$ws.Range("E20").Value = "2026-02-12 22:49:32"
$ws.Range("O20").Value = "-3.6 °C"
$ws.Range("E21").Value = "2026-02-12 22:49:35"
$ws.Range("J21").Value = "1002.9 hPa"
$ws.Range("E22").Value = "2026-02-12 22:49:38"
$ws.Range("E23").Value = "2026-02-12 22:49:41"
$ws.Range("E24").Value = "2026-02-12 22:49:44"
$ws.Range("J24").Value = "1006.9 hPa"
$ws.Range("E25").Value = "2026-02-12 22:49:46"
$ws.Range("E26").Value = "2026-02-12 22:49:49"
$ws.Range("J26").Value = "999.6 hPa"
$ws.Range("N26").Value = "2.0 °C 22:26 TU"
$ws.Range("O26").Value = "5.7 °C"
$ws.Range("E27").Value = "2026-02-12 22:49:52"
$ws.Range("E28").Value = "2026-02-12 22:49:55"
$ws.Range("J28").Value = "999.7 hPa"
$ws.Range("O28").Value = "13.6 °C"
$ws.Range("E29").Value = "2026-02-12 22:49:57"
$ws.Range("O29").Value = "13.8 °C"
$ws.Range("E30").Value = "2026-02-12 22:50:00"
$ws.Range("E31").Value = "2026-02-12 22:50:03"
$ws.Range("J31").Value = "999.4 hPa"
$ws.Range("O31").Value = "14.1 °C"
$ws.Range("E32").Value = "2026-02-12 22:50:05"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "58%"
$ws.Range("E33").Value = "2026-02-12 22:50:08"
$ws.Range("J33").Value = "1002.1 hPa"
$ws.Range("N33").Value = "2.7 °C 22:29 TU"
$ws.Range("O33").Value = "6.4 °C"
$ws.Range("E34").Value = "2026-02-12 22:50:11"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "58%"
$ws.Range("E35").Value = "2026-02-12 22:50:13"
$ws.Range("E36").Value = "2026-02-12 22:50:16"
$ws.Range("J36").Value = "1000.3 hPa"
$ws.Range("O36").Value = "14.3 °C"
$ws.Range("E37").Value = "2026-02-12 22:50:19"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "50%"
$ws.Range("J37").Value = "1001.1 hPa"
$ws.Range("N37").Value = "2.6 °C 22:21 TU"
$ws.Range("O37").Value = "9.5 °C"
$ws.Range("E38").Value = "2026-02-12 22:50:22"
$ws.Range("N38").Value = "12.3 °C 22:29 TU"
$ws.Range("O38").Value = "15.7 °C"
$ws.Range("E39").Value = "2026-02-12 22:50:24"
$ws.Range("E40").Value = "2026-02-12 22:50:27"
$ws.Range("J40").Value = "1003.7 hPa"
$ws.Range("N40").Value = "3.6 °C 22:29 TU"
$ws.Range("O40").Value = "9.2 °C"
$ws.Range("E41").Value = "2026-02-12 22:50:30"
$ws.Range("E42").Value = "2026-02-12 22:50:33"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "63%"
$ws.Range("N42").Value = "7.5 °C 22:10 TU"
$ws.Range("O42").Value = "13.7 °C"
$ws.Range("E43").Value = "2026-02-12 22:50:35"
$ws.Range("E44").Value = "2026-02-12 22:50:38"
$ws.Range("E45").Value = "2026-02-12 22:50:41"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "55%"
$ws.Range("J45").Value = "1005.5 hPa"
$ws.Range("N45").Value = "1.9 °C 22:24 TU"
$ws.Range("O45").Value = "6.7 °C"
$ws.Range("E46").Value = "2026-02-12 22:50:44"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "41%"
$ws.Range("O46").Value = "15.7 °C"
